$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(463).Insert()

$ws.Range("A463").Value = 4
$ws.Range("B463").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C463").Value = "Los Lagos"
$ws.Range("D463").Value = 45218
$ws.Range("E463").Value = 10
$ws.Range("F463").Value = 100112040
$ws.Range("G463").Value = "Cilantro"
$ws.Range("H463").Value = "Sin especificar"
$ws.Range("I463").Value = "Primera"
$ws.Range("J463").Value = 80
$ws.Range("K463").Value = 12000
$ws.Range("L463").Value = 12000
$ws.Range("M463").Value = 12000
$ws.Range("N463").Value = "$/caja 36 atados"
$ws.Range("O463").Value = "Región Metropolitana"
$ws.Range("P463").Value = 333
$ws.Range("Q463").Value = 36
$ws.Range("R463").Value = "Hortaliza"
